# chore: update Sheets via scheduled runner
#
# Refreshes the market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) on each job sheet's Leve table with newly scraped values.
# Only cell values change -- no structural/formatting edits.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2999.4285
$ws.Range("J70").Value = 3439.2
$ws.Range("L70").Value = 10317.6
$ws.Range("N70").Value = -10857.6
$ws.Range("H73").Value = 2999.4285
$ws.Range("J73").Value = 3439.2
$ws.Range("L73").Value = 10317.6
$ws.Range("N73").Value = -12189.6
$ws.Range("H74").Value = 4872.5
$ws.Range("I74").Value = 4250
$ws.Range("K74").Value = 4250
$ws.Range("M74").Value = -3314
$ws.Range("H77").Value = 4872.5
$ws.Range("I77").Value = 4250
$ws.Range("K77").Value = 21250
$ws.Range("M77").Value = -16570
$ws.Range("H80").Value = 34528.332
$ws.Range("I80").Value = 86635.57000000001
$ws.Range("K80").Value = 259906.71
$ws.Range("M80").Value = -258908.71
$ws.Range("H83").Value = 34528.332
$ws.Range("I83").Value = 86635.57000000001
$ws.Range("K83").Value = 779720.1300000001
$ws.Range("M83").Value = -774728.1300000001
$ws.Range("H86").Value = 71545.39999999999
$ws.Range("I86").Value = 76298.64
$ws.Range("K86").Value = 76298.64
$ws.Range("M86").Value = -75175.64
$ws.Range("H89").Value = 71545.39999999999
$ws.Range("I89").Value = 76298.64
$ws.Range("K89").Value = 381493.2
$ws.Range("M89").Value = -375877.2
$ws.Range("H92").Value = 66949.734
$ws.Range("I92").Value = 83628.414
$ws.Range("K92").Value = 83628.414
$ws.Range("M92").Value = -82380.414
$ws.Range("H111").Value = 1579.875
$ws.Range("J111").Value = 2999.6667
$ws.Range("L111").Value = 8999.000100000001
$ws.Range("N111").Value = -15133.0001
$ws.Range("H125").Value = 1049.1428
$ws.Range("J125").Value = 1154.0834
$ws.Range("L125").Value = 10386.7506
$ws.Range("N125").Value = -15306.7506
$ws.Range("H131").Value = 3804.7222
$ws.Range("J131").Value = 5732.5
$ws.Range("L131").Value = 17197.5
$ws.Range("N131").Value = -27277.5
$ws.Range("H132").Value = 4386.6963
$ws.Range("I132").Value = 2735.102
$ws.Range("K132").Value = 8205.306
$ws.Range("M132").Value = -5675.306
$ws.Range("H137").Value = 3227.0938
$ws.Range("I137").Value = 885.5238000000001
$ws.Range("J137").Value = 7697.364
$ws.Range("K137").Value = 2656.5714
$ws.Range("L137").Value = 23092.092
$ws.Range("M137").Value = -106.5714000000003
$ws.Range("N137").Value = -28192.092
$ws.Range("H141").Value = 9402.210999999999
$ws.Range("I141").Value = 8867.177
$ws.Range("K141").Value = 26601.531
$ws.Range("M141").Value = -21421.531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11820.787
$ws.Range("I32").Value = 8188.3057
$ws.Range("J32").Value = 23708.908
$ws.Range("K32").Value = 8188.3057
$ws.Range("L32").Value = 23708.908
$ws.Range("M32").Value = -7901.3057
$ws.Range("N32").Value = -24282.908
$ws.Range("H61").Value = 3548
$ws.Range("I61").Value = 1897.2858
$ws.Range("K61").Value = 1897.2858
$ws.Range("M61").Value = -1685.2858
$ws.Range("H102").Value = 1800.875
$ws.Range("I102").Value = 1622.2069
$ws.Range("K102").Value = 1622.2069
$ws.Range("M102").Value = -0.206899999999905
$ws.Range("H132").Value = 39984.258
$ws.Range("I132").Value = 41387.5
$ws.Range("K132").Value = 124162.5
$ws.Range("M132").Value = -121632.5
$ws.Range("H136").Value = 3548
$ws.Range("I136").Value = 1897.2858
$ws.Range("K136").Value = 5691.857400000001
$ws.Range("M136").Value = -3141.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1016.41174
$ws.Range("I94").Value = 585.8333
$ws.Range("K94").Value = 585.8333
$ws.Range("M94").Value = -134.8333
$ws.Range("H107").Value = 1186.4419
$ws.Range("I107").Value = 898
$ws.Range("J107").Value = 2276.111
$ws.Range("K107").Value = 898
$ws.Range("L107").Value = 2276.111
$ws.Range("M107").Value = 1022
$ws.Range("N107").Value = -6116.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4175.6284
$ws.Range("I31").Value = 2965.72
$ws.Range("J31").Value = 7200.4
$ws.Range("K31").Value = 2965.72
$ws.Range("L31").Value = 7200.4
$ws.Range("M31").Value = -2670.72
$ws.Range("N31").Value = -7790.4
$ws.Range("H34").Value = 4175.6284
$ws.Range("I34").Value = 2965.72
$ws.Range("J34").Value = 7200.4
$ws.Range("K34").Value = 2965.72
$ws.Range("L34").Value = 7200.4
$ws.Range("M34").Value = -2763.72
$ws.Range("N34").Value = -7604.4
$ws.Range("H105").Value = 3919.5715
$ws.Range("I105").Value = 984.5
$ws.Range("K105").Value = 984.5
$ws.Range("M105").Value = 762.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 960.9091
$ws.Range("I8").Value = 960.9091
$ws.Range("K8").Value = 2882.7273
$ws.Range("M8").Value = -2743.7273
$ws.Range("H18").Value = 887.3333
$ws.Range("I18").Value = 887.3333
$ws.Range("K18").Value = 2661.9999
$ws.Range("M18").Value = -2492.9999
$ws.Range("H40").Value = 270.6
$ws.Range("I40").Value = 284.66666
$ws.Range("K40").Value = 1138.66664
$ws.Range("M40").Value = -1069.66664
$ws.Range("H86").Value = 647.1667
$ws.Range("I86").Value = 577.6667
$ws.Range("J86").Value = 716.6667
$ws.Range("K86").Value = 1733.0001
$ws.Range("L86").Value = 2150.0001
$ws.Range("M86").Value = -547.0001
$ws.Range("N86").Value = -4522.0001
$ws.Range("H89").Value = 647.1667
$ws.Range("I89").Value = 577.6667
$ws.Range("J89").Value = 716.6667
$ws.Range("K89").Value = 5199.0003
$ws.Range("L89").Value = 6450.0003
$ws.Range("M89").Value = 728.9997000000003
$ws.Range("N89").Value = -18306.0003
$ws.Range("H121").Value = 961.1
$ws.Range("I121").Value = 466.33334
$ws.Range("J121").Value = 1173.1428
$ws.Range("K121").Value = 1399.00002
$ws.Range("L121").Value = 3519.4284
$ws.Range("M121").Value = -89.00001999999995
$ws.Range("N121").Value = -6139.428400000001
$ws.Range("H122").Value = 753.7
$ws.Range("I122").Value = 676.8570999999999
$ws.Range("J122").Value = 933
$ws.Range("K122").Value = 6091.7139
$ws.Range("L122").Value = 8397
$ws.Range("M122").Value = -3641.7139
$ws.Range("N122").Value = -13297
$ws.Range("H123").Value = 5166.3335
$ws.Range("I123").Value = 2999.75
$ws.Range("J123").Value = 9499.5
$ws.Range("K123").Value = 8999.25
$ws.Range("L123").Value = 28498.5
$ws.Range("M123").Value = -6549.25
$ws.Range("N123").Value = -33398.5
$ws.Range("H129").Value = 449195.78
$ws.Range("I129").Value = 1545.6
$ws.Range("J129").Value = 609070.9
$ws.Range("K129").Value = 4636.799999999999
$ws.Range("L129").Value = 1827212.7
$ws.Range("M129").Value = 363.2000000000007
$ws.Range("N129").Value = -1837212.7
$ws.Range("H131").Value = 4775570
$ws.Range("I131").Value = 2918
$ws.Range("J131").Value = 6267023.5
$ws.Range("K131").Value = 8754
$ws.Range("L131").Value = 18801070.5
$ws.Range("M131").Value = -3714
$ws.Range("N131").Value = -18811150.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1580.2106
$ws.Range("I97").Value = 1468.4
$ws.Range("K97").Value = 1468.4
$ws.Range("M97").Value = -972.4000000000001
$ws.Range("H102").Value = 1846.3928
$ws.Range("I102").Value = 1321.0952
$ws.Range("J102").Value = 3422.2856
$ws.Range("K102").Value = 1321.0952
$ws.Range("L102").Value = 3422.2856
$ws.Range("M102").Value = 300.9048
$ws.Range("N102").Value = -6666.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 92526.55
$ws.Range("I22").Value = 100779.3
$ws.Range("J22").Value = 9999
$ws.Range("K22").Value = 100779.3
$ws.Range("L22").Value = 9999
$ws.Range("M22").Value = -100484.3
$ws.Range("N22").Value = -10589
$ws.Range("H27").Value = 92526.55
$ws.Range("I27").Value = 100779.3
$ws.Range("J27").Value = 9999
$ws.Range("K27").Value = 100779.3
$ws.Range("L27").Value = 9999
$ws.Range("M27").Value = -100672.3
$ws.Range("N27").Value = -10213
$ws.Range("H46").Value = 3764.4736
$ws.Range("J46").Value = 4807.75
$ws.Range("L46").Value = 4807.75
$ws.Range("N46").Value = -5183.75
$ws.Range("H61").Value = 1943.5862
$ws.Range("I61").Value = 1956.3846
$ws.Range("J61").Value = 1832.6666
$ws.Range("K61").Value = 1956.3846
$ws.Range("L61").Value = 1832.6666
$ws.Range("M61").Value = -1754.3846
$ws.Range("N61").Value = -2236.6666
$ws.Range("H68").Value = 5629.8335
$ws.Range("I68").Value = 3999
$ws.Range("J68").Value = 5956
$ws.Range("K68").Value = 3999
$ws.Range("L68").Value = 5956
$ws.Range("M68").Value = -3250
$ws.Range("N68").Value = -7454
$ws.Range("H71").Value = 5629.8335
$ws.Range("I71").Value = 3999
$ws.Range("J71").Value = 5956
$ws.Range("K71").Value = 19995
$ws.Range("L71").Value = 29780
$ws.Range("M71").Value = -16251
$ws.Range("N71").Value = -37268
$ws.Range("H82").Value = 2378
$ws.Range("I82").Value = 1188.4667
$ws.Range("J82").Value = 3270.15
$ws.Range("K82").Value = 1188.4667
$ws.Range("L82").Value = 3270.15
$ws.Range("M82").Value = -827.4666999999999
$ws.Range("N82").Value = -3992.15
$ws.Range("H85").Value = 2378
$ws.Range("I85").Value = 1188.4667
$ws.Range("J85").Value = 3270.15
$ws.Range("K85").Value = 1188.4667
$ws.Range("L85").Value = 3270.15
$ws.Range("M85").Value = 59.53330000000005
$ws.Range("N85").Value = -5766.15
$ws.Range("H93").Value = 1609.9231
$ws.Range("I93").Value = 1492.2858
$ws.Range("K93").Value = 1492.2858
$ws.Range("M93").Value = -244.2858000000001
$ws.Range("H113").Value = 1943.5862
$ws.Range("I113").Value = 1956.3846
$ws.Range("J113").Value = 1832.6666
$ws.Range("K113").Value = 1956.3846
$ws.Range("L113").Value = 1832.6666
$ws.Range("M113").Value = 213.6153999999999
$ws.Range("N113").Value = -6172.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1344.7778
$ws.Range("I100").Value = 1125.1111
$ws.Range("J100").Value = 1784.1111
$ws.Range("K100").Value = 2250.2222
$ws.Range("L100").Value = 2250.2222
$ws.Range("M100").Value = -1709.2222
$ws.Range("N100").Value = -4650.2222
